$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the _GoBack bookmark from its original (now-orphaned)
# location -- the paragraph that used to hold only the bookmark becomes a
# plain empty paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: delete the run containing the "Using ChatGPT to refine..."
# paragraph text (the whole run is removed, leaving the rest of the
# paragraph -- which ends with a single trailing space run -- intact).
# ---------------------------------------------------------------------------
$target = "Using ChatGPT to refine certain sections of the report proved to be beneficial. By providing ChatGPT with specific sections and requesting revisions, it was able to enhance clarity, precision, and grammatical accuracy. Iterative refinement often involved multiple iterations to ensure that the resulting paragraphs maintained the original meaning, sounded professional, and avoided unnecessary complexity. This collaborative process facilitated the creation of polished and concise sections within the report."
$findRange = $d.Content
$found = $findRange.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $findRange.Delete()
}

# ---------------------------------------------------------------------------
# Change 3: insert four new empty paragraphs (matching the "ind left=360"
# formatting of the existing signature-block paragraphs) right before the
# last of those paragraphs, and put the _GoBack bookmark back on the new
# 4th (last) of the inserted paragraphs.
# ---------------------------------------------------------------------------

# Locate the paragraph immediately preceding the trailing "ind left=360"
# paragraph that sits right before the honor-code statement, by searching
# for that statement and walking back to the previous paragraph.
$honorRange = $d.Content
$honorFound = $honorRange.Find.Execute("I have neither given nor received unauthorized aid", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$honorParaIndex = $honorRange.Paragraphs(1).Range.Start
$lastIndentPara = $d.Range($honorRange.Start, $honorRange.Start).Paragraphs(1).Previous(1)

$insertionPoint = $d.Range($lastIndentPara.Range.End, $lastIndentPara.Range.End)

$newParasXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:ind w:left="360"/></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:ind w:left="360"/></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:ind w:left="360"/></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:ind w:left="360"/></w:pPr></w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newParasXml)

# The 4th newly-inserted paragraph (the one that should carry the bookmark)
# is the paragraph right before the original trailing "ind left=360"
# paragraph (which followed our insertion point).
$bookmarkPara = $lastIndentPara.Next(4)
$bookmarkPara.Range.Select()
$d.Bookmarks.Add("_GoBack", $word.Selection.Range)

Write-Output "done"
